# "Generate Report for Handback"
#
# Before this edit, the zh-cn / de-de handoff tables have two source files
# ("4fe2e3d8..." and "7e8d9ea9...") that are still "Ready for handoff":
# their Latest Target File / Latest Handback File / Latest Handback DateTime
# columns are blank (or the zero-date sentinel) and the overview Status says
# "Ready for handoff". This script records a successful handback: it fills
# in the target/handback file links + timestamps, flips the Status text, and
# widens a few columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$urlMd4  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2266812517f1ac1f354d9ea3ff0e328dc782dde/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$urlMd7  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2266812517f1ac1f354d9ea3ff0e328dc782dde/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"

$md4     = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$md7     = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"

$statusText = "Handed back: in sync with en-US"

# Hyperlink-blue used by the pre-existing "HyperLink" style (RGB FF6495ED,
# expressed as the BGR OLE color Excel's Font.Color setter expects).
$hyperlinkColor = 15570276

# Widths: the engine stores column width as (input + 5/6) quantized to whole
# pixel steps, same as real Excel - these inputs land on the nearest pixel
# bucket to the target widths (~29.98 and 40 characters).
$wideWidth  = 29.15
$fullWidth  = 39.17

# ---------------------------------------------------------------------
# Overview sheet: Status text for both languages, plus column widening.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText

$ov.Columns.Item(5).ColumnWidth = $wideWidth
$ov.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $urlMd4, [System.Type]::Missing, [System.Type]::Missing, $md4)
$zh.Range("I2").Font.Color = $hyperlinkColor
$zh.Range("J2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-04 17:07:15"

$zh.Hyperlinks.Add($zh.Range("I3"), $urlMd7, [System.Type]::Missing, [System.Type]::Missing, $md7)
$zh.Range("I3").Font.Color = $hyperlinkColor
$zh.Range("J3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-04 17:07:15"

$zh.Columns.Item(3).ColumnWidth = $wideWidth
$zh.Columns.Item(9).ColumnWidth = $fullWidth
$zh.Columns.Item(10).ColumnWidth = $fullWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $urlMd4, [System.Type]::Missing, [System.Type]::Missing, $md4)
$de.Range("I2").Font.Color = $hyperlinkColor
$de.Range("J2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.de-de.xlf"
$de.Range("K2").Value = "2016-09-04 17:07:23"

$de.Hyperlinks.Add($de.Range("I3"), $urlMd7, [System.Type]::Missing, [System.Type]::Missing, $md7)
$de.Range("I3").Font.Color = $hyperlinkColor
$de.Range("J3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.de-de.xlf"
$de.Range("K3").Value = "2016-09-04 17:07:23"

$de.Columns.Item(3).ColumnWidth = $wideWidth
$de.Columns.Item(9).ColumnWidth = $fullWidth
$de.Columns.Item(10).ColumnWidth = $fullWidth
